# CellValues.xlsx fixture update
#
# ClosedXML's XLCell now loads the *file's cached value* into
# XLCell.CachedValue instead of re-deriving it. For a TimeSpan-backed
# cell the cached "GetFormattedString()" result therefore changes from
# the naive ToString() rendering ("1.02:31:45") to the value Excel
# actually produced using the cell's real number format ([h]:mm:ss),
# which is "26:31:45" for the elapsed-time value stored in row 7.
#
# Sheet "Cell Values": G7 ("GetFormattedString()" column) gets the
# corrected cached text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cell Values")
$ws.Range("G7").Value2 = "26:31:45"

# Column C ("Using Get...()") also ends up a bit narrower in the
# refreshed fixture.
$ws.Columns.Item(3).ColumnWidth = 8.95
